$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (Coke Oven Coke) keeps its fuel label, but its subfuel code gets
# a real code instead of the placeholder "x".
$ws.Range("C6").Value = "02_01_coke_oven_coke"

# Insert four new coal/oil rows directly below the (now updated) Coke Oven
# Coke row: BKB/PB, Peat, Peat Products, Crude Oil.
$ws.Range("A7:A10").EntireRow.Insert()
$ws.Range("A7").Value = "BKB/PB"
$ws.Range("B7").Value = "02_coal_products"
$ws.Range("C7").Value = "02_08_bkb_pb"

$ws.Range("A8").Value = "Peat"
$ws.Range("B8").Value = "03_peat"
$ws.Range("C8").Value = "x"

$ws.Range("A9").Value = "Peat  Products"
$ws.Range("B9").Value = "04_peat_products"
$ws.Range("C9").Value = "x"

$ws.Range("A10").Value = "Crude Oil"
$ws.Range("B10").Value = "06_crude_oil_and_ngl"
$ws.Range("C10").Value = "06_01_crude_oil"

# Rows 11-19 (previously 7-15: Motor Gasoline ... Geothermal) are unchanged
# and simply shifted down by the insert above.

# Insert two new biomass rows below Geothermal (row 19): Charcoal, Other biomass.
$ws.Range("A20:A21").EntireRow.Insert()
$ws.Range("A20").Value = "Charcoal"
$ws.Range("B20").Value = "15_solid_biomass"
$ws.Range("C20").Value = "15_03_charcoal"

$ws.Range("A21").Value = "Other biomass"
$ws.Range("B21").Value = "15_solid_biomass"
$ws.Range("C21").Value = "15_05_other_biomass"

# Rows 22-23 (previously 16-17: Other Biomass, Biogas) are unchanged and
# simply shifted down.

# Insert one new row below Biogas (row 23): Industrial Waste.
$ws.Range("A24").EntireRow.Insert()
$ws.Range("A24").Value = "Industrial Waste"
$ws.Range("B24").Value = "16_others"
$ws.Range("C24").Value = "16_02_industrial_waste"

# Rows 25-29 (previously 18-22: Municipal Solid Waste (Renewable/Non-renewable),
# Biodiesel, Electricity, Heat) are unchanged and simply shifted down.

# Update the saved selection to match the author's final cursor position.
$ws.Range("G23").Select()
